$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.728.25"
$ws.Range("E2").Value = "  -1.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.099.74"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.85"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5192"
$ws.Range("E7").Value = "  +0.63%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4390"
$ws.Range("E8").Value = "  -1.46%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.92"
$ws.Range("E9").Value = "  +2.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09246"
$ws.Range("E10").Value = "  +3.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.169"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.72"
$ws.Range("E12").Value = "  -4.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.102.56"
$ws.Range("E13").Value = "  -0.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.801"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.98"
$ws.Range("E16").Value = "  +3.88%  "

$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.01"
$ws.Range("E19").Value = "  +0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06668"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.220"
$ws.Range("E22").Value = "  -0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.779.38"
$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.56"
$ws.Range("E24").Value = "  -2.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.307"
$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.348.87"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.19"
$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.496"
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.73"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("E31").Value = "  -4.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.703"
$ws.Range("E32").Value = "  +4.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.200"
$ws.Range("E34").Value = "  -0.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.951"
$ws.Range("E35").Value = "  -0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.379"
$ws.Range("E36").Value = "  +7.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.43"
$ws.Range("E37").Value = "  +0.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02580"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06720"
$ws.Range("E39").Value = "  -1.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7003"
$ws.Range("E40").Value = "  +2.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.53"
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.333"
$ws.Range("E42").Value = "  +4.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2217"
$ws.Range("E43").Value = "  -4.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6819"
$ws.Range("E44").Value = "  +6.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.41"
$ws.Range("E45").Value = "  +0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.344"
$ws.Range("E46").Value = "  +1.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000359"
$ws.Range("E47").Value = "  -2.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.622"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.198"
$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.61"
$ws.Range("E51").Value = "  -1.52%  "
